$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 131  # was 130
$ws.Range("F4").Value = 1268  # was 1266
$ws.Range("F7").Value = 969  # was 966
$ws.Range("F11").Value = 99  # was 98
$ws.Range("F12").Value = 653  # was 654
$ws.Range("F13").Value = 913  # was 907
$ws.Range("F14").Value = 1797  # was 1793
$ws.Range("F15").Value = 3850  # was 3813
$ws.Range("F16").Value = 1144  # was 1137
$ws.Range("F17").Value = 109  # was 106
$ws.Range("F18").Value = 2567  # was 2556
$ws.Range("F20").Value = 1070  # was 1065
$ws.Range("F21").Value = 3547  # was 3530
$ws.Range("F22").Value = 744  # was 738
$ws.Range("F23").Value = 836  # was 833
$ws.Range("F25").Value = 2202  # was 2199
$ws.Range("F26").Value = 107  # was 106
$ws.Range("F27").Value = 825  # was 823
$ws.Range("F28").Value = 164  # was 162
$ws.Range("F29").Value = 284  # was 251
$ws.Range("F30").Value = 190  # was 187
$ws.Range("F32").Value = 1319  # was 1309
$ws.Range("F33").Value = 1929  # was 1923
$ws.Range("F34").Value = 484  # was 479
$ws.Range("F35").Value = 25  # was 11
$ws.Range("F37").Value = 585  # was 582
$ws.Range("F38").Value = 269  # was 265
$ws.Range("F39").Value = 10  # was 4
$ws.Range("F40").Value = 165  # was 164
$ws.Range("F41").Value = 228  # was 226
$ws.Range("F42").Value = 74  # was 73

$ws = $wb.Worksheets.Item(2)
$ws.Range("F12").Value = 112  # was 109

$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 414  # was 400

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 414  # was 400
$ws.Range("F3").Value = 131  # was 130
$ws.Range("F4").Value = 1268  # was 1266
$ws.Range("F6").Value = 969  # was 966
$ws.Range("F13").Value = 99  # was 98
$ws.Range("F15").Value = 913  # was 907
$ws.Range("F16").Value = 1797  # was 1793
$ws.Range("F17").Value = 3850  # was 3813
$ws.Range("F18").Value = 1144  # was 1137
$ws.Range("F19").Value = 109  # was 106
$ws.Range("F21").Value = 2568  # was 2556
$ws.Range("F23").Value = 1070  # was 1065
$ws.Range("F24").Value = 3547  # was 3530
$ws.Range("F25").Value = 744  # was 738
$ws.Range("F26").Value = 836  # was 833
$ws.Range("F29").Value = 2202  # was 2199
$ws.Range("F33").Value = 107  # was 106
$ws.Range("F34").Value = 112  # was 109
$ws.Range("F35").Value = 825  # was 823
$ws.Range("F36").Value = 164  # was 162
$ws.Range("F37").Value = 285  # was 251
$ws.Range("F38").Value = 190  # was 187
$ws.Range("F41").Value = 1319  # was 1309
$ws.Range("F42").Value = 1929  # was 1923
$ws.Range("F44").Value = 484  # was 479
$ws.Range("F45").Value = 585  # was 582
$ws.Range("F46").Value = 270  # was 265
$ws.Range("F47").Value = 165  # was 164
$ws.Range("F48").Value = 228  # was 226
$ws.Range("F49").Value = 74  # was 73
